$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 30's phone number ("A30") currently holds the phone number as text;
# duplicate the whole row to become the new row 31 first so the copy keeps
# the original text-typed phone number, then correct A30 to a real number
# and update the new row's varying fields.
$ws.Range("A30:I30").Copy($ws.Range("A31:I31"))

# Row 30: phone number should be stored as a number.
$ws.Range("A30").Value = 71277620

# Row 31: new payment record for 71277620 (Cash) at 2025-08-18T17:04:26.
$ws.Range("D31").Value = "2025-08-18T17:04:26"
$ws.Range("E31").Value = 760
$ws.Range("G31").Value = 0
$ws.Range("H31").Value = 0
$ws.Range("I31").Value = 760
